# PSP Time Recording Log - add two new logged activities (10/14, 10/15)
# to rows 18 and 19 of the author sheet ("작성자명"), matching a "Node JS 강의"
# activity already used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 18: 10월 14일, 18:00-21:00, 0 min interruption, 180 min delta ---
$ws.Range("A18").Value = "10월 14일"
$ws.Range("A18").Characters(3, 5).Font.Name = "Arial Unicode MS"
$ws.Range("A18").Characters(3, 5).Font.Size = 10

$ws.Range("B18").Value = 0.75
$ws.Range("C18").Value = 0.875
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 180
$ws.Range("F18").Value = "Node JS 강의"

# --- Row 19: 10월 15일, 18:00-22:30, 60 min interruption, 210 min delta ---
$ws.Range("A19").Value = "10월 15일"
$ws.Range("A19").Characters(3, 5).Font.Name = "Arial Unicode MS"
$ws.Range("A19").Characters(3, 5).Font.Size = 10

$ws.Range("B19").Value = 0.75
$ws.Range("C19").Value = 0.9375
$ws.Range("D19").Value = 60
$ws.Range("E19").Value = 210
$ws.Range("F19").Value = "Node JS 강의"

# Match the author's final selection/cursor position on the sheet.
$ws.Range("F20").Select()
